$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Reuse the existing header style/formatting (bold, bordered, centered) by
# copying the format from an existing header cell (E1) onto the new ones.
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# New boolean (FALSE) data cells F2:H4
$ws.Range("F2:H4").Value = $false
